$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 69.40000000000001
$ws.Range("I5").Value = 63.75
$ws.Range("K5").Value = 63.75
$ws.Range("M5").Value = 51.25
$ws.Range("H12").Value = 525225
$ws.Range("I12").Value = 400
$ws.Range("J12").Value = 700166.7
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 700166.7
$ws.Range("M12").Value = -230
$ws.Range("N12").Value = -700506.7
$ws.Range("H106").Value = 673.7
$ws.Range("I106").Value = 673.7
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 673.7
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -42.70000000000005
$ws.Range("H107").Value = 1289.0588
$ws.Range("I107").Value = 1201
$ws.Range("K107").Value = 1201
$ws.Range("M107").Value = 719
$ws.Range("H121").Value = 1169.125
$ws.Range("J121").Value = 1913.25
$ws.Range("L121").Value = 5739.75
$ws.Range("N121").Value = -9233.75
$ws.Range("H132").Value = 3126868.8
$ws.Range("I132").Value = 3450002
$ws.Range("J132").Value = 3246.6667
$ws.Range("K132").Value = 10350006
$ws.Range("L132").Value = 9740.000100000001
$ws.Range("M132").Value = -10347476
$ws.Range("N132").Value = -14800.0001
$ws.Range("H138").Value = 5007.9565
$ws.Range("I138").Value = 2205.3845
$ws.Range("J138").Value = 6112
$ws.Range("K138").Value = 6616.1535
$ws.Range("L138").Value = 18336
$ws.Range("M138").Value = -1476.1535
$ws.Range("N138").Value = -28616

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 53.75
$ws.Range("I5").Value = 55
$ws.Range("K5").Value = 55
$ws.Range("M5").Value = 57
$ws.Range("H63").Value = 2995.8333
$ws.Range("I63").Value = 3195
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 3195
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -2509
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 2995.8333
$ws.Range("I66").Value = 3195
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 15975
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -12543
$ws.Range("N66").Value = -16864
$ws.Range("H74").Value = 887.8333
$ws.Range("I74").Value = 813.4211
$ws.Range("J74").Value = 1016.36365
$ws.Range("K74").Value = 813.4211
$ws.Range("L74").Value = 1016.36365
$ws.Range("M74").Value = 60.57889999999998
$ws.Range("N74").Value = -2764.36365
$ws.Range("H77").Value = 887.8333
$ws.Range("I77").Value = 813.4211
$ws.Range("J77").Value = 1016.36365
$ws.Range("K77").Value = 4067.1055
$ws.Range("L77").Value = 5081.81825
$ws.Range("M77").Value = 300.8944999999999
$ws.Range("N77").Value = -13817.81825
$ws.Range("H110").Value = 1219.5
$ws.Range("I110").Value = 700.6875
$ws.Range("K110").Value = 700.6875
$ws.Range("M110").Value = 1344.3125
$ws.Range("H132").Value = 16950956
$ws.Range("I132").Value = 24391500
$ws.Range("J132").Value = 3046.1667
$ws.Range("K132").Value = 73174500
$ws.Range("L132").Value = 9138.500100000001
$ws.Range("M132").Value = -73171970
$ws.Range("N132").Value = -14198.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 53.75
$ws.Range("I4").Value = 55
$ws.Range("K4").Value = 55
$ws.Range("M4").Value = 60
$ws.Range("H22").Value = 420
$ws.Range("H94").Value = 684.2432
$ws.Range("I94").Value = 622.9167
$ws.Range("J94").Value = 797.46155
$ws.Range("K94").Value = 622.9167
$ws.Range("L94").Value = 797.46155
$ws.Range("M94").Value = -171.9167
$ws.Range("N94").Value = -1699.46155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1432.6
$ws.Range("I22").Value = 530
$ws.Range("J22").Value = 1658.25
$ws.Range("K22").Value = 530
$ws.Range("L22").Value = 1658.25
$ws.Range("M22").Value = -180
$ws.Range("N22").Value = -2358.25
$ws.Range("H31").Value = 2684.6667
$ws.Range("I31").Value = 1587.1538
$ws.Range("K31").Value = 1587.1538
$ws.Range("M31").Value = -1292.1538
$ws.Range("H34").Value = 2684.6667
$ws.Range("I34").Value = 1587.1538
$ws.Range("K34").Value = 1587.1538
$ws.Range("M34").Value = -1385.1538
$ws.Range("H58").Value = 8476771
$ws.Range("I58").Value = 1404.5526
$ws.Range("J58").Value = 23813148
$ws.Range("K58").Value = 1404.5526
$ws.Range("L58").Value = 23813148
$ws.Range("M58").Value = -1201.5526
$ws.Range("N58").Value = -23813554
$ws.Range("H132").Value = 2160.9
$ws.Range("I132").Value = 1586.1212
$ws.Range("J132").Value = 4870.5713
$ws.Range("K132").Value = 4758.363600000001
$ws.Range("L132").Value = 14611.7139
$ws.Range("M132").Value = -2228.363600000001
$ws.Range("N132").Value = -19671.7139
$ws.Range("H134").Value = 1171.4902
$ws.Range("J134").Value = 3381.0908
$ws.Range("L134").Value = 10143.2724
$ws.Range("N134").Value = -15213.2724
$ws.Range("H136").Value = 8476771
$ws.Range("I136").Value = 1404.5526
$ws.Range("J136").Value = 23813148
$ws.Range("K136").Value = 4213.6578
$ws.Range("L136").Value = 71439444
$ws.Range("M136").Value = -1663.6578
$ws.Range("N136").Value = -71444544

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2880.75
$ws.Range("I80").Value = 2505
$ws.Range("J80").Value = 3006
$ws.Range("K80").Value = 2505
$ws.Range("L80").Value = 3006
$ws.Range("M80").Value = -1507
$ws.Range("N80").Value = -5002
$ws.Range("H83").Value = 2880.75
$ws.Range("I83").Value = 2505
$ws.Range("J83").Value = 3006
$ws.Range("K83").Value = 12525
$ws.Range("L83").Value = 15030
$ws.Range("M83").Value = -7533
$ws.Range("N83").Value = -25014
$ws.Range("H97").Value = 2425.0527
$ws.Range("I97").Value = 1424
$ws.Range("J97").Value = 4594
$ws.Range("K97").Value = 1424
$ws.Range("L97").Value = 4594
$ws.Range("M97").Value = -928
$ws.Range("N97").Value = -5586
$ws.Range("H102").Value = 49876.76
$ws.Range("I102").Value = 1299.9333
$ws.Range("J102").Value = 171318.83
$ws.Range("K102").Value = 1299.9333
$ws.Range("L102").Value = 171318.83
$ws.Range("M102").Value = 322.0667000000001
$ws.Range("N102").Value = -174562.83
$ws.Range("H126").Value = 2961.1365
$ws.Range("I126").Value = 1876.8182
$ws.Range("J126").Value = 4045.4546
$ws.Range("K126").Value = 5630.4546
$ws.Range("L126").Value = 12136.3638
$ws.Range("M126").Value = -3160.4546
$ws.Range("N126").Value = -17076.3638
$ws.Range("H132").Value = 3434
$ws.Range("I132").Value = 3529.3572
$ws.Range("J132").Value = 3363.7368
$ws.Range("K132").Value = 10588.0716
$ws.Range("L132").Value = 10091.2104
$ws.Range("M132").Value = -8058.071599999999
$ws.Range("N132").Value = -15151.2104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3171.2258
$ws.Range("I132").Value = 2227.0908
$ws.Range("J132").Value = 3690.5
$ws.Range("K132").Value = 6681.2724
$ws.Range("L132").Value = 11071.5
$ws.Range("M132").Value = -4151.2724
$ws.Range("N132").Value = -16131.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4799.857
$ws.Range("I62").Value = 4533
$ws.Range("K62").Value = 4533
$ws.Range("M62").Value = -3909
$ws.Range("H65").Value = 4799.857
$ws.Range("I65").Value = 4533
$ws.Range("K65").Value = 22665
$ws.Range("M65").Value = -19545
